# Update RMSE-Z analysis output values (IDPT v1.0 release data) for rows 2-5,
# columns C:X. Columns A and B are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1057.366153846154,0.9449212416697137,1.592307692307702,104.845008475194,1.543827478235422,-0.1827674016766958,0.5463759282705529,-0.04236414324339088,0.07486626928564244,0.2338355087905244,-0.4793858822651518,52.69538461538462,453.9027468816307,384.9541656106993,0,0.3953038831653787,0.1324823337841791,0.4871414100856032,-1.543827478235422,104.845008475194,0.8991878264002933,325),
    @(1044.558631921824,0.9273634173733174,1.207817589576557,224.065873356555,0.8052034913387703,-0.5034705838343624,0.9135237439530184,-0.04550794939280208,0.03620005125635661,0.4082240437299039,-0.1805969095931438,48.92508143322475,431.1477472113035,396.3819937731997,0,0.3964857532455751,0.09663089433481362,0.6025802970168854,-0.8052034913387703,224.065873356555,1.269836697463687,614),
    @(1072.215988779804,0.9255930217044702,2.608274894810668,355.9308675936459,1.941059737938881,-0.7624430437987492,1.471880691620011,0.05637119804769756,-0.04982625946089946,0.381716158358394,0.1291709516199683,36.31977559607293,332.0157809135414,360.5822558223263,0,0.4792164035118819,-0.005253785763080949,0.6044449404197529,-1.941059737938881,355.9308675936459,1.755233503250185,713),
    @(1083.763076923077,0.9250313410370359,3.192307692307701,453.9033896499955,2.286065233174964,-1.054429674149614,1.931815550366022,0.004669425456218405,-0.03394537228809891,0.4915647704794313,0.5371939252963939,21.09846153846154,205.4677907829984,275.8228375760975,0,0.4468266442576123,0.01237566954056361,0.6897165914182772,-2.286065233174964,453.9033896499955,2.293092874221573,325)
)

# Columns C (3) through X (24)
$startCol = 3
$startRow = 2

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($startRow + $r, $startCol + $c).Value = $rowVals[$c]
    }
}
